$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# Row 12 (TestScenario_3 / Edit Case): mark as Rejected
$ws.Range("I12").Value = "Rejected"

# Row 20 (TestScenario_4 / Delete Case): mark as Rejected, reason "Testworked"
$ws.Range("I20").Value = "Rejected"
$ws.Range("J20").Value = "Testworked"

# Row 12 reason: "Tested"
$ws.Range("J12").Value = "Tested"

# Update selection to match the saved view state (I12 was last active cell)
$ws.Range("I12").Select()
